$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update inserts two new daily price records (date serial 44516,
# "Lane Late" variety) for Femacal de La Calera - Naranja, pushing all the
# existing rows from 509 downward by two rows.
$ws.Rows("509:510").Insert()

# New row 509: Lane Late / Primera
$ws.Range("A509").Value = 3
$ws.Range("B509").Value = "Femacal de La Calera"
$ws.Range("C509").Value = "Coquimbo"
$ws.Range("D509").Value = 44516
$ws.Range("E509").Value = 5
$ws.Range("F509").Value = "Fruta"
$ws.Range("G509").Value = 100102
$ws.Range("H509").Value = "Cítricos"
$ws.Range("I509").Value = 100102005
$ws.Range("J509").Value = "Naranja"
$ws.Range("K509").Value = "Lane Late"
$ws.Range("L509").Value = "Primera"
$ws.Range("M509").Value = 163
$ws.Range("N509").Value = 5000
$ws.Range("O509").Value = 6000
$ws.Range("P509").Value = 5521
$ws.Range("Q509").Value = "`$/malla 13 kilos"
$ws.Range("R509").Value = "Provincia de Quillota"
$ws.Range("S509").Value = 425
$ws.Range("T509").Value = 13

# New row 510: Lane Late / Segunda
$ws.Range("A510").Value = 3
$ws.Range("B510").Value = "Femacal de La Calera"
$ws.Range("C510").Value = "Coquimbo"
$ws.Range("D510").Value = 44516
$ws.Range("E510").Value = 5
$ws.Range("F510").Value = "Fruta"
$ws.Range("G510").Value = 100102
$ws.Range("H510").Value = "Cítricos"
$ws.Range("I510").Value = 100102005
$ws.Range("J510").Value = "Naranja"
$ws.Range("K510").Value = "Lane Late"
$ws.Range("L510").Value = "Segunda"
$ws.Range("M510").Value = 150
$ws.Range("N510").Value = 4000
$ws.Range("O510").Value = 4500
$ws.Range("P510").Value = 4267
$ws.Range("Q510").Value = "`$/malla 13 kilos"
$ws.Range("R510").Value = "Provincia de Quillota"
$ws.Range("S510").Value = 328
$ws.Range("T510").Value = 13
